$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PAYC")

# Row 4 (Inventory)
$ws.Range("B4").Value = 1000000.0
$ws.Range("C4").Value = 1000000.0
$ws.Range("D4").Value = 1000000.0
$ws.Range("E4").Value = 1000000.0
$ws.Range("F4").Value = 1000000.0

# Row 13 (Accounts Payable)
$ws.Range("B13").Value = 7000000.0
$ws.Range("C13").Value = 9000000.0
$ws.Range("D13").Value = 4000000.0
$ws.Range("E13").Value = 6000000.0
$ws.Range("F13").Value = 5000000.0

# Row 23 (Long Term Tax Liability (Deferred))
$ws.Range("B23").Value = 113000000.0
$ws.Range("C23").Value = 102000000.0
$ws.Range("D23").Value = 99000000.0
$ws.Range("E23").Value = 93000000.0
$ws.Range("F23").Value = 91000000.0
